$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 67) with the next quarterly observation.
$row = 67

$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 1).Value = (Get-Date -Year 2025 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0).Date

$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 0.5
$ws.Cells.Item($row, 4).Value = 2.5
$ws.Cells.Item($row, 5).Value = 3.5
$ws.Cells.Item($row, 6).Value = 6.5
$ws.Cells.Item($row, 7).Value = 12.5
$ws.Cells.Item($row, 8).Value = 15.5
